$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The generated report lost its "Docentes responsaveis" detail rows (the two
# unlabeled rows holding the professors' names) as well as the long free-text
# bodies for "Objetivos", "Programa resumido", "Programa" and "Bibliografia".
# What's left behind cascades the remaining short answers up by one slot in
# each of those sections. Reproduce that exact end state:
#
#   1) Each "long text" row's B/C content is overwritten by the text that
#      used to sit one section below it (a stray leftover value).
#   2) The two anonymous rows carrying the professors' names are removed
#      entirely, shifting every following row up by two.
#
# Use Range.Copy so the destination cells pick up the exact same shared
# string + style as their source (a plain .Value assignment would have
# Excel "helpfully" reinterpret date-shaped text like "01/01/2019" as a
# real date serial, which is not what happened here).

# -- Cascade the trailing long-text rows, deepest dependency first --------
$ws.Range("B22").Copy($ws.Range("B23"))   # Bibliografia:          <- old Norma de recuperacao text
$ws.Range("C22").Copy($ws.Range("C23"))

$ws.Range("B21").Copy($ws.Range("B22"))   # Norma de recuperação:  <- old Critério text
$ws.Range("C21").Copy($ws.Range("C22"))

$ws.Range("B20").Copy($ws.Range("B21"))   # Critério:              <- old Método text
$ws.Range("C20").Copy($ws.Range("C21"))

$ws.Range("B8").Copy($ws.Range("B20"))    # Método:                <- old Ativação text (01/01/2019)
$ws.Range("C8").Copy($ws.Range("C20"))

# -- Objetivos / Programa resumido / Programa pick up the stray docente names
$ws.Range("B13").Copy($ws.Range("B10"))   # Objetivos:             <- "6634418 - Antonio Clelio Ribeiro"
$ws.Range("C13").Copy($ws.Range("C10"))

$ws.Range("B13").Copy($ws.Range("B15"))   # Programa resumido:     <- "6634418 - Antonio Clelio Ribeiro"
$ws.Range("C13").Copy($ws.Range("C15"))

$ws.Range("B14").Copy($ws.Range("B17"))   # Programa:              <- "1285870 - Marcos Villela Barcza"
$ws.Range("C14").Copy($ws.Range("C17"))

# -- Remove the two now-redundant "docente name" rows ----------------------
$ws.Rows("13:14").Delete()
